# ---------------------------------------------------------------------------
# "Refined metadata to be additional tab"
#
# 1) Refresh the "time_taken" timestamps in column F of the existing "data"
#    sheet (rows 2-73) to reflect the latest panel query run.
# 2) Add a new "metadata" worksheet (placed right after "data") carrying the
#    panel-level query metadata that used to live only implicitly.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$times = @(
    "2021-10-05 14:22:33.979375",
    "2021-10-05 14:22:33.979383",
    "2021-10-05 14:22:33.979386",
    "2021-10-05 14:22:33.979388",
    "2021-10-05 14:22:33.979391",
    "2021-10-05 14:22:33.979403",
    "2021-10-05 14:22:33.979405",
    "2021-10-05 14:22:33.979408",
    "2021-10-05 14:22:33.979410",
    "2021-10-05 14:22:33.979413",
    "2021-10-05 14:22:33.979415",
    "2021-10-05 14:22:33.979418",
    "2021-10-05 14:22:33.979420",
    "2021-10-05 14:22:33.979423",
    "2021-10-05 14:22:33.979425",
    "2021-10-05 14:22:33.979428",
    "2021-10-05 14:22:33.979430",
    "2021-10-05 14:22:33.979433",
    "2021-10-05 14:22:33.979458",
    "2021-10-05 14:22:33.979472",
    "2021-10-05 14:22:33.979475",
    "2021-10-05 14:22:33.979478",
    "2021-10-05 14:22:33.979481",
    "2021-10-05 14:22:33.979483",
    "2021-10-05 14:22:33.979486",
    "2021-10-05 14:22:33.979489",
    "2021-10-05 14:22:33.979491",
    "2021-10-05 14:22:33.979494",
    "2021-10-05 14:22:33.979496",
    "2021-10-05 14:22:33.979498",
    "2021-10-05 14:22:33.979501",
    "2021-10-05 14:22:33.979503",
    "2021-10-05 14:22:33.979506",
    "2021-10-05 14:22:33.979509",
    "2021-10-05 14:22:33.979511",
    "2021-10-05 14:22:33.979514",
    "2021-10-05 14:22:33.979516",
    "2021-10-05 14:22:33.979519",
    "2021-10-05 14:22:33.979521",
    "2021-10-05 14:22:33.979523",
    "2021-10-05 14:22:33.979526",
    "2021-10-05 14:22:33.979529",
    "2021-10-05 14:22:33.979531",
    "2021-10-05 14:22:33.979534",
    "2021-10-05 14:22:33.979536",
    "2021-10-05 14:22:33.979538",
    "2021-10-05 14:22:33.979541",
    "2021-10-05 14:22:33.979543",
    "2021-10-05 14:22:33.979546",
    "2021-10-05 14:22:33.979548",
    "2021-10-05 14:22:33.979551",
    "2021-10-05 14:22:33.979553",
    "2021-10-05 14:22:33.979556",
    "2021-10-05 14:22:33.979558",
    "2021-10-05 14:22:33.979561",
    "2021-10-05 14:22:33.979563",
    "2021-10-05 14:22:33.979566",
    "2021-10-05 14:22:33.979568",
    "2021-10-05 14:22:33.979571",
    "2021-10-05 14:22:33.979573",
    "2021-10-05 14:22:33.979576",
    "2021-10-05 14:22:33.979578",
    "2021-10-05 14:22:33.979580",
    "2021-10-05 14:22:33.979583",
    "2021-10-05 14:22:33.979587",
    "2021-10-05 14:22:33.979590",
    "2021-10-05 14:22:33.979593",
    "2021-10-05 14:22:33.979595",
    "2021-10-05 14:22:33.979598",
    "2021-10-05 14:22:33.979600",
    "2021-10-05 14:22:33.979603",
    "2021-10-05 14:22:33.979605"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $times[$i]
}

# --- add the new "metadata" tab, right after "data" -----------------------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Respiratory ciliopathies including non-CF bronchiectasis"
$meta.Range("C2").Value = 550
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.45"
$meta.Range("E2").Value = "2021-05-04T10:21:45.830038Z"
$meta.Range("F2").Value = "2021-10-05 14:22:33.975924"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/550/?format=json"

# Match the header / index-column styling used on the "data" sheet (bold
# text, thin box border, centered + top-aligned) for B1:G1 and A2.
$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexCell = $meta.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.Borders.LineStyle = 1
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160

# Keep "data" as the active tab (unchanged from before the edit)
$ws.Activate()
